# Mod Part Patching for 0.8.0
# Fill in the missing Tier value for the rtg-0625 part (row 4) so its
# generated ModuleManager patch formula resolves "reactors" -> "reactors8",
# matching the pattern already used by the surrounding rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = 8

# Leave the selection on F7:F8, matching where the edit was made.
$ws.Range("F7:F8").Select() | Out-Null
